$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.25517891189561936
$ws.Cells.Item(2, 1).Value = -0.0059999999642599278
$ws.Cells.Item(3, 1).Value = -0.0039999999664210861
$ws.Cells.Item(4, 1).Value = -0.0079999999392317278
$ws.Cells.Item(5, 1).Value = -0.0029999999627481344
$ws.Cells.Item(6, 1).Value = -0.0019999999553270698
$ws.Cells.Item(7, 1).Value = -0.009999999909534818
$ws.Cells.Item(8, 1).Value = -0.02130194861916701
$ws.Cells.Item(9, 1).Value = -0.0019999999538682367
$ws.Cells.Item(10, 1).Value = -0.0019999999544442204
$ws.Cells.Item(11, 1).Value = -0.0029999999488099505
$ws.Cells.Item(12, 1).Value = -0.0034999999463907194
$ws.Cells.Item(13, 1).Value = -0.0034999999486133859
$ws.Cells.Item(14, 1).Value = -0.0079999999242224007
$ws.Cells.Item(15, 1).Value = 0.03669175600679786
$ws.Cells.Item(16, 1).Value = -0.0019999999589486173
$ws.Cells.Item(17, 1).Value = -0.0019999999580253558
$ws.Cells.Item(18, 1).Value = -0.0039999999465223368
$ws.Cells.Item(19, 1).Value = -0.047642278871539823
$ws.Cells.Item(20, 1).Value = -0.0039999999728905777
$ws.Cells.Item(21, 1).Value = -0.0039999999725948143
$ws.Cells.Item(22, 1).Value = -0.0039999999723647761
$ws.Cells.Item(23, 1).Value = -0.0049999999589225297
$ws.Cells.Item(24, 1).Value = -0.019999999870788265
$ws.Cells.Item(25, 1).Value = -0.019999999869092733
$ws.Cells.Item(26, 1).Value = -0.0024999999505528336
$ws.Cells.Item(27, 1).Value = -0.0024999999480139756
$ws.Cells.Item(28, 1).Value = -0.001999999939398478
$ws.Cells.Item(29, 1).Value = -0.0069999999031482574
$ws.Cells.Item(30, 1).Value = -0.059999999600523157
$ws.Cells.Item(31, 1).Value = -0.0069999998961840504
$ws.Cells.Item(32, 1).Value = -0.0099999998788220523
$ws.Cells.Item(33, 1).Value = -0.0039999999124464836
